# Insert a new weekly observation row into the daily price table.
# A new row is inserted at row 40 (pushing all following rows down by
# one, extending the used range from A1:R116 to A1:R117) and populated
# with the new "Arveja Verde" price observation for 2023-11-08.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 40; Excel shifts rows
# 40..116 down to 41..117 and copies row 40's formatting (date style
# on column D) down into the new row.
$ws.Rows.Item(40).Insert()

$newRow = 40

$ws.Cells.Item($newRow, 1).Value  = 7
$ws.Cells.Item($newRow, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($newRow, 3).Value  = "Ñuble"
$ws.Cells.Item($newRow, 4).Value  = Get-Date -Year 2023 -Month 11 -Day 8 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item($newRow, 5).Value  = 16
$ws.Cells.Item($newRow, 6).Value  = 100112022
$ws.Cells.Item($newRow, 7).Value  = "Arveja Verde"
$ws.Cells.Item($newRow, 8).Value  = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value  = "Primera"
$ws.Cells.Item($newRow, 10).Value = 30
$ws.Cells.Item($newRow, 11).Value = 23000
$ws.Cells.Item($newRow, 12).Value = 23000
$ws.Cells.Item($newRow, 13).Value = 23000
$ws.Cells.Item($newRow, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item($newRow, 15).Value = "Región del Maule"
$ws.Cells.Item($newRow, 16).Value = 920
$ws.Cells.Item($newRow, 17).Value = 25
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
